$d = $word.ActiveDocument

# --- Paragraph 1: "Foi criado em 1991 e a primeira versão lançada em 94." ---
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Foi criado em 1991 e a primeira versão lançada em 94.</w:t></w:r></w:p>
'@

# --- Paragraph 2: empty paragraph (no numbering), indent left=1080 ---
$p2xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>
'@

# --- Paragraph 3: "Aula 2 – O Que As Versões Do Linux Tem Em Comum:" ---
$p3xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Aula 2 – O Que As Versões Do Linux Tem Em Comum:</w:t></w:r></w:p>
'@

# --- Paragraph 4: " " + "Linux é licenciado como GPL, ..." ---
$p4xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Linux é licenciado como GPL, o que permite que todos possam alterar tudo e criar suas próprias versões do Linux.</w:t></w:r></w:p>
'@

# The document currently ends with the paragraph "Linus teve uma sacada ..."
# (the last paragraph before the sectPr). Each new paragraph is created by
# splitting off a fresh paragraph after the previous last one, then that
# fresh paragraph's content/formatting is replaced wholesale via InsertXML
# so the resulting <w:p> matches the target exactly (no leftover numbering,
# rsid, or paraId attributes).

$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$d.Paragraphs($lastIndex).Range.InsertXML($p1xml)

$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$d.Paragraphs($lastIndex).Range.InsertXML($p2xml)

$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$d.Paragraphs($lastIndex).Range.InsertXML($p3xml)

$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
$lastIndex = $lastIndex + 1
$d.Paragraphs($lastIndex).Range.InsertXML($p4xml)
